# Fix(all scripts): Updating routes for __tmp__ files
# Update the selected "Formatos" filter value from the generic
# "--todos--" placeholder to the specific "Supermercado" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A2").Value = "Supermercado"
